$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Cells whose type/style changes (number <-> "N/A" text) ---
$ws.Range("D15").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").Formula = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("J14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 6

$ws.Range("J14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 2

$ws.Range("K14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 200

$ws.Range("G22").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").Formula = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("C26").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Formula = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("F28").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$ws.Range("F29").Formula = "'0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("J14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

$ws.Range("J14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

$ws.Range("J14").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("I30").Value = 1

$excel.CutCopyMode = $false

# --- Plain value updates (style/type unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = -50
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 19
$ws.Range("H16").Value = 46.153846153846
$ws.Range("I16").Value = 66
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 20
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -16.455696202531
$ws.Range("N16").Value = -68.571428571428
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 63.157894736842
$ws.Range("I17").Value = 97
$ws.Range("J17").Value = 88
$ws.Range("K17").Value = 10.227272727272
$ws.Range("L17").Value = 25.974025974026
$ws.Range("M17").Value = 42.647058823529
$ws.Range("N17").Value = 42.647058823529
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 242.857142857143
$ws.Range("I18").Value = 64
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 72.972972972973
$ws.Range("L18").Value = 45.454545454545
$ws.Range("M18").Value = -25.581395348837
$ws.Range("N18").Value = -85.454545454545
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 14
$ws.Range("I19").Value = 161
$ws.Range("J19").Value = 192
$ws.Range("K19").Value = -16.145833333333
$ws.Range("L19").Value = 62.626262626262
$ws.Range("M19").Value = 93.975903614457
$ws.Range("N19").Value = 9.523809523809
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -37.5
$ws.Range("F20").Value = 28
$ws.Range("H20").Value = -6.666666666666
$ws.Range("I20").Value = 113
$ws.Range("J20").Value = 108
$ws.Range("K20").Value = 4.629629629629
$ws.Range("L20").Value = 88.333333333333
$ws.Range("M20").Value = 85.245901639344
$ws.Range("N20").Value = -75.802997858672
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -6.451612903225
$ws.Range("F21").Value = 160
$ws.Range("G21").Value = 120
$ws.Range("H21").Value = 33.333333333333
$ws.Range("I21").Value = 505
$ws.Range("J21").Value = 487
$ws.Range("K21").Value = 3.696098562628
$ws.Range("L21").Value = 51.651651651651
$ws.Range("M21").Value = 31.510416666666
$ws.Range("N21").Value = -62.369597615499
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 75
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 16.666666666666
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 44.444444444444
$ws.Range("I23").Value = 36
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = 44
$ws.Range("L23").Value = 56.521739130434
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 36.666666666666
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = 3.883495145631
$ws.Range("I24").Value = 408
$ws.Range("J24").Value = 318
$ws.Range("K24").Value = 28.301886792452
$ws.Range("L24").Value = 51.111111111111
$ws.Range("M24").Value = 76.623376623376
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 23.076923076923
$ws.Range("I25").Value = 150
$ws.Range("J25").Value = 125
$ws.Range("K25").Value = 20
$ws.Range("L25").Value = 54.639175257732
$ws.Range("M25").Value = -11.242603550295
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = -18.181818181818
$ws.Range("D27").Value = 4
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 6.666666666666
$ws.Range("L27").Value = 33.333333333333
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -63.636363636363
$ws.Range("N28").Value = -69.230769230769
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -55.555555555555
$ws.Range("N29").Value = -60
